# The workbook originally has a 3-column / 3-row table starting at A1:
#   Row1: A,B,C (headers)
#   Row2: alpha, beta, gamma
#   Row3: delta, epsilon, zeta
#
# The edit reflows the sheet to match the standard upload template: a new
# instructional line is added at A1, the header row moves down to row 4,
# and the data rows move down to rows 10 and 11 (leaving blank spacer rows
# as the template dictates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the whole original table down by 3 rows, so the header row (was
# row 1) lands on row 4.
$ws.Rows("1:3").Insert()

# Push the data rows (now rows 5 and 6) further down so they land on the
# template's rows 10 and 11, leaving blank rows 5-9 in between.
$ws.Rows("5:9").Insert()

# Add the new instructional text in the now-empty A1.
$ws.Range("A1").Value2 = "The location of the data in this spreadsheet matches the upload template:"

# The previous selection (D4) no longer corresponds to anything meaningful
# after the reflow; return the selection to the default top-left cell.
$ws.Range("A1").Select()

$wb.Save()
